$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (columns B-F, I-N); column G stays 1, column H stays blank.
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.022585203763145, 1.027219167439481, 1.048508784654992, 1.052778160187862, 1.030627091185019, 1.027770213808853, 1.030039534561403, 1.051268432637377, 1.05552595220033, 1.029229765399353)
    3 = @(1.02, 1.023485738831831, 1.02785232462888, 1.049580120962158, 1.053870428705247, 1.030758492797049, 1.028309115442415, 1.030480842930338, 1.052151157720566, 1.056430417696365, 1.029769432335048)
    4 = @(1.02, 1.024068792292468, 1.028262134072192, 1.050274417906803, 1.05457810382439, 1.030842232359159, 1.028657557496632, 1.030765825835823, 1.052722812070963, 1.057015976797008, 1.030118369216949)
    5 = @(1.02, 1.024313990070817, 1.028434443945593, 1.050566556061499, 1.054875826267003, 1.030877128080782, 1.028803978599067, 1.030885494486349, 1.052963248037031, 1.057262219255952, 1.030264998254143)
    6 = @(1.02, 1.024355164642835, 1.028463377014311, 1.050615622320971, 1.05492582781144, 1.030882969127225, 1.028828559571173, 1.030905579249554, 1.053003624873039, 1.057303568733977, 1.03028961413405)
    7 = @(1.02, 1.024072068315796, 1.028264436384827, 1.050278320467079, 1.054582081159049, 1.030842699849889, 1.028659514234919, 1.030767425398277, 1.052726024346854, 1.057019266813165, 1.030120328734028)
    8 = @(1.02, 1.022889471354928, 1.027433121344162, 1.048870625707457, 1.053147109784616, 1.030671764961139, 1.027952392127775, 1.03018879488374, 1.051566655598448, 1.055831556073411, 1.029412202432384)
    9 = @(1.02, 1.020808285731431, 1.02596917614774, 1.046398315084925, 1.050625463141739, 1.030360728313628, 1.026704375147539, 1.029164826763844, 1.049527341204195, 1.053741057466194, 1.028162413124886)
    10 = @(1.02, 1.019422701921524, 1.024993926386973, 1.044755674271672, 1.048949081025129, 1.030146792709158, 1.025871086615385, 1.028479314820933, 1.04817027878975, 1.05234904319975, 1.027327941227437)
    11 = @(1.02, 1.018823184751063, 1.024571819089003, 1.04404571963154, 1.04822431510098, 1.030052602190634, 1.025509971332724, 1.028181812360894, 1.047583250410113, 1.051746684903279, 1.026966313119655)
    12 = @(1.02, 1.018600565833014, 1.024415058409073, 1.043782209608019, 1.047955273360202, 1.030017382483674, 1.025375793424698, 1.028071206762153, 1.04736529080603, 1.051523001807268, 1.026831944663609)
    13 = @(1.02, 1.018648315234837, 1.024448682783283, 1.043838724401457, 1.048012976088973, 1.030024947773632, 1.025404577001644, 1.028094936546363, 1.047412039852105, 1.051570979885111, 1.026860769116535)
    14 = @(1.02, 1.018804781600332, 1.024558860607542, 1.044023933753777, 1.048202072584368, 1.030049695672258, 1.025498881029862, 1.028172671701052, 1.047565231985671, 1.051728193969906, 1.02695520706729)
    15 = @(1.02, 1.018901194789192, 1.024626748654204, 1.044138073697786, 1.048318603532636, 1.030064912788075, 1.025556979099846, 1.028220553679823, 1.047659630505094, 1.051825066669519, 1.027013387643198)
    16 = @(1.02, 1.019462499452756, 1.025021944216442, 1.044802819467313, 1.048997204965795, 1.030153011101493, 1.025895046485423, 1.028499045014421, 1.048209250368802, 1.052389028063626, 1.027351935123239)
    17 = @(1.02, 1.019814712073838, 1.025269889728711, 1.04522015057319, 1.049423173573586, 1.030207856859085, 1.026107028606692, 1.028673556197827, 1.048554170284824, 1.052742891762289, 1.027564218283435)
    18 = @(1.02, 1.020020195073096, 1.025414529646607, 1.045463700043589, 1.049671741747491, 1.030239697396343, 1.026230645590423, 1.028775280790209, 1.048755412997724, 1.052949332359451, 1.027688010817462)
    19 = @(1.02, 1.02009026678953, 1.025463851046323, 1.04554676569654, 1.049756515355276, 1.03025052872891, 1.026272790910683, 1.028809955223655, 1.048824041176582, 1.053019729678359, 1.027730215988911)
    20 = @(1.02, 1.019776918504017, 1.02524328570216, 1.045175361724691, 1.049377460009594, 1.030201987946933, 1.026084287883556, 1.02865483949831, 1.048517157767765, 1.052704921613387, 1.027541445265864)
    21 = @(1.02, 1.0187587042588, 1.024526415179301, 1.043969388725684, 1.048146383721223, 1.030042414465405, 1.025471112041257, 1.02814978338698, 1.04752011826362, 1.051681896690894, 1.026927398643535)
    22 = @(1.02, 1.018118909033923, 1.024075857333829, 1.043212295727385, 1.04737333289896, 1.029940735554707, 1.025085332132055, 1.027831656799626, 1.046893754373235, 1.05103902539047, 1.026541070882614)
    23 = @(1.02, 1.018458038710893, 1.024314690273441, 1.043613535874991, 1.047783049100162, 1.029994765157523, 1.02528986496151, 1.028000356184117, 1.047225752812515, 1.051379790739304, 1.026745894172147)
    24 = @(1.02, 1.019793995657557, 1.02525530686269, 1.045195599484051, 1.049398115668042, 1.030204640319426, 1.026094563525916, 1.028663296969562, 1.048533881955623, 1.052722078572176, 1.027551735500814)
    25 = @(1.02, 1.021345996111537, 1.026347521338687, 1.047036488105084, 1.051276540928386, 1.03044230030454, 1.027027246400553, 1.029430056430717, 1.050054118025146, 1.054281213197449, 1.028485742892113)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}
